# Update cryptos list figures (prices / 1h volume change %) to the latest
# scrape, matching the GitHub Actions "Updated cryptos list" commit.
# Cells in column D sometimes look like plain decimals (e.g. "113.05"), so we
# force NumberFormat "@" (Text) before assigning them to stop Excel from
# auto-converting them to numeric values; values with two dots (e.g.
# "43.733.15") are never auto-converted and are left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.733.15'
$ws.Range('E2').Value = '  -0.02%  '

$ws.Range('D3').Value = '2.293.75'
$ws.Range('E3').Value = '  -0.12%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '113.05'
$ws.Range('E5').Value = '  +16.37%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '269.31'
$ws.Range('E6').Value = '  +0.21%  '

$ws.Range('E7').Value = '  +0.60%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.618'
$ws.Range('E9').Value = '  +1.20%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '48.16'
$ws.Range('E10').Value = '  +5.91%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0951'
$ws.Range('E11').Value = '  +1.75%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '9.11'
$ws.Range('E12').Value = '  +14.86%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.106'
$ws.Range('E13').Value = '  +0.61%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '15.78'
$ws.Range('E14').Value = '  +0.26%  '

$ws.Range('D15').Value = '2.636.95'
$ws.Range('E15').Value = '  -0.15%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.851'
$ws.Range('E16').Value = '  -0.89%  '

$ws.Range('D17').Value = '2.290.49'
$ws.Range('E17').Value = '  -0.36%  '

$ws.Range('D18').Value = '43.703.79'
$ws.Range('E18').Value = '  -0.09%  '

$ws.Range('E19').Value = '  -0.67%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.74'
$ws.Range('E20').Value = '  +8.90%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '72.32'
$ws.Range('E21').Value = '  +0.36%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.46'
$ws.Range('E22').Value = '  -2.43%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '232.48'
$ws.Range('E23').Value = '  -0.36%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.81'
$ws.Range('E24').Value = '  +7.45%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.81'
$ws.Range('E25').Value = '  +6.62%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.68'
$ws.Range('E27').Value = '  +3.33%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '42.04'
$ws.Range('E28').Value = '  +7.97%  '

$ws.Range('E29').Value = '  -1.96%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.26'
$ws.Range('E30').Value = '  +1.46%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '175.39'
$ws.Range('E31').Value = '  +0.12%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '21.52'
$ws.Range('E32').Value = '  -1.92%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0926'
$ws.Range('E33').Value = '  +2.55%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.69'
$ws.Range('E34').Value = '  +5.02%  '

$ws.Range('E35').Value = '  +1.32%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.66'
$ws.Range('E36').Value = '  +2.95%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0362'
$ws.Range('E37').Value = '  +2.77%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.108'
$ws.Range('E38').Value = '  +0.07%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.84'
$ws.Range('E39').Value = '  +14.02%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '13.87'
$ws.Range('E40').Value = '  +13.61%  '

$ws.Range('B41').Value = 'MultiversX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '74.15'
$ws.Range('E41').Value = '  +15.81%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.243'
$ws.Range('E42').Value = '  +1.15%  '

$ws.Range('B43').Value = 'LidoDAOToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.39'
$ws.Range('E43').Value = '  +1.96%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '6.31'
$ws.Range('E44').Value = '  +22.15%  '

$ws.Range('E45').Value = '  +0.19%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.39'
$ws.Range('E46').Value = '  +2.51%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.78'
$ws.Range('E47').Value = '  -0.37%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '102.94'
$ws.Range('E48').Value = '  +5.54%  '

$ws.Range('E49').Value = '  -1.96%  '

$ws.Range('E50').Value = '  +2.87%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.461'
$ws.Range('E51').Value = '  +6.12%  '
